$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 29 ("「芸術」..." entry) entirely; all rows below shift up by one.
$ws.Rows.Item(29).Delete()
